$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = ""

$ws.Range("H62").Value = 24000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 24000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = ""
$ws.Range("M62").Value = 24000
$ws.Range("N62").Value = -25248

$ws.Range("H65").Value = 24000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 24000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = ""
$ws.Range("M65").Value = 120000
$ws.Range("N65").Value = -126240

$ws.Range("H98").Value = 738.75
$ws.Range("I98").Value = 587.6
$ws.Range("K98").Value = 587.6
$ws.Range("M98").Value = 910.4

$ws.Range("H113").Value = 2997.5
$ws.Range("I113").Value = 3995
$ws.Range("K113").Value = 3995
$ws.Range("M113").Value = -741

$ws.Range("H122").Value = 738.75
$ws.Range("I122").Value = 587.6
$ws.Range("K122").Value = 1762.8
$ws.Range("M122").Value = 687.1999999999998

$ws.Range("H132").Value = 1880.9286
$ws.Range("I132").Value = 2036.5
$ws.Range("K132").Value = 6109.5
$ws.Range("M132").Value = -3579.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 6100
$ws.Range("I31").Value = 6100
$ws.Range("K31").Value = 6100
$ws.Range("M31").Value = -5806

$ws.Range("H32").Value = 12221.083
$ws.Range("I32").Value = 7962.1665
$ws.Range("K32").Value = 7962.1665
$ws.Range("M32").Value = -7675.1665

$ws.Range("H97").Value = 664.8421
$ws.Range("I97").Value = 454.375
$ws.Range("J97").Value = 1787.3334
$ws.Range("K97").Value = 454.375
$ws.Range("L97").Value = 1787.3334
$ws.Range("M97").Value = 41.625
$ws.Range("N97").Value = -2779.3334

$ws.Range("H122").Value = 324861.34
$ws.Range("I122").Value = 501310.1
$ws.Range("K122").Value = 1503930.3
$ws.Range("M122").Value = -1501480.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3450.875
$ws.Range("I99").Value = 3823.8
$ws.Range("J99").Value = 3184.5
$ws.Range("K99").Value = 3823.8
$ws.Range("L99").Value = 3184.5
$ws.Range("M99").Value = -2325.8
$ws.Range("N99").Value = -6180.5

$ws.Range("H105").Value = 5430.857
$ws.Range("I105").Value = 3666.6365
$ws.Range("J105").Value = 11899.667
$ws.Range("K105").Value = 3666.6365
$ws.Range("L105").Value = 11899.667
$ws.Range("M105").Value = -1919.6365
$ws.Range("N105").Value = -15393.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 299.66666
$ws.Range("I22").Value = 299
$ws.Range("K22").Value = 299
$ws.Range("M22").Value = 51

$ws.Range("H62").Value = 70698.164
$ws.Range("I62").Value = 4838
$ws.Range("K62").Value = 4838
$ws.Range("M62").Value = -4214

$ws.Range("H65").Value = 70698.164
$ws.Range("I65").Value = 4838
$ws.Range("K65").Value = 24190
$ws.Range("M65").Value = -21070

$ws.Range("H134").Value = 3081.8333
$ws.Range("I134").Value = 2345.6155
$ws.Range("K134").Value = 7036.8465
$ws.Range("M134").Value = -4501.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 711.1111
$ws.Range("J9").Value = 799.875
$ws.Range("L9").Value = 2399.625
$ws.Range("N9").Value = -2847.625

$ws.Range("H12").Value = 867.5454999999999
$ws.Range("I12").Value = 899.5
$ws.Range("J12").Value = 860.44446
$ws.Range("K12").Value = 2698.5
$ws.Range("L12").Value = 2581.33338
$ws.Range("M12").Value = -2525.5
$ws.Range("N12").Value = -2927.33338

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = ""

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = ""

$ws.Range("H36").Value = 1883.3334
$ws.Range("I36").Value = 1883.3334
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 5650.0002
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = ""
$ws.Range("N36").Value = -5481.0002

$ws.Range("H75").Value = 336.8
$ws.Range("J75").Value = 296
$ws.Range("L75").Value = 888
$ws.Range("N75").Value = -2884

$ws.Range("H78").Value = 336.8
$ws.Range("J78").Value = 296
$ws.Range("L78").Value = 2664
$ws.Range("N78").Value = -12648

$ws.Range("H111").Value = 525.75
$ws.Range("I111").Value = 525.75
$ws.Range("K111").Value = 1577.25
$ws.Range("M111").Value = 1489.75

$ws.Range("H112").Value = 3649.8333
$ws.Range("I112").Value = 950
$ws.Range("J112").Value = 4999.75
$ws.Range("K112").Value = 2850
$ws.Range("L112").Value = 14999.25
$ws.Range("M112").Value = -1742
$ws.Range("N112").Value = -17215.25

$ws.Range("H117").Value = 2776.625
$ws.Range("I117").Value = 750
$ws.Range("K117").Value = 2250
$ws.Range("M117").Value = 1192

$ws.Range("H120").Value = 15375
$ws.Range("J120").Value = 15954.546
$ws.Range("L120").Value = 47863.638
$ws.Range("N120").Value = -57539.638

$ws.Range("H137").Value = 5685.273
$ws.Range("I137").Value = 3788.1428
$ws.Range("K137").Value = 11364.4284
$ws.Range("M137").Value = -6264.428400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4467.1
$ws.Range("I16").Value = 4467.1
$ws.Range("K16").Value = 4467.1
$ws.Range("M16").Value = -4297.1

$ws.Range("H22").Value = 8256.714
$ws.Range("I22").Value = 1549
$ws.Range("K22").Value = 1549
$ws.Range("M22").Value = -1254

$ws.Range("H27").Value = 8256.714
$ws.Range("I27").Value = 1549
$ws.Range("K27").Value = 1549
$ws.Range("M27").Value = -1442

$ws.Range("H46").Value = 2895.5186
$ws.Range("J46").Value = 3825.3333
$ws.Range("L46").Value = 3825.3333
$ws.Range("N46").Value = -4201.3333

$ws.Range("H55").Value = 245.90909
$ws.Range("I55").Value = 255.71428
$ws.Range("K55").Value = 255.71428
$ws.Range("M55").Value = -82.71428

$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251

$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256

$ws.Range("H93").Value = 799
$ws.Range("I93").Value = 799
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 799
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = ""
$ws.Range("N93").Value = 449

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 58599.8
$ws.Range("J64").Value = 65666.664
$ws.Range("L64").Value = 65666.664
$ws.Range("N64").Value = -66162.664

$ws.Range("H67").Value = 58599.8
$ws.Range("J67").Value = 65666.664
$ws.Range("L67").Value = 65666.664
$ws.Range("N67").Value = -67382.664

$ws.Range("H76").Value = 48748.75
$ws.Range("J76").Value = 48331.668
$ws.Range("L76").Value = 48331.668
$ws.Range("N76").Value = -48961.668

$ws.Range("H79").Value = 48748.75
$ws.Range("J79").Value = 48331.668
$ws.Range("L79").Value = 48331.668
$ws.Range("N79").Value = -50515.668

$ws.Range("H126").Value = 3173.2856
$ws.Range("I126").Value = 2552
$ws.Range("K126").Value = 7656
$ws.Range("M126").Value = -5186

